$d = $word.ActiveDocument

# Locate the paragraph that ends with "Learn math in English." (the last
# paragraph of the body, right before the sectPr).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Learn math in English.*") {
        $targetPara = $p
    }
}

$insertionPoint = $targetPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

# Re-fetch paragraphs by index: the newly created paragraph is the one right
# after $targetPara.
$idx = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $targetPara.Range.Start) {
        $idx = $i
        break
    }
}
$newPara = $d.Paragraphs.Item($idx + 1)
$newRange = $newPara.Range

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p>' + `
      '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
      '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">3, 2^3  is read as </w:t></w:r>' + `
      '<w:r><w:t>' + $openQuote + '</w:t></w:r>' + `
      '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>two to the powe of three</w:t></w:r>' + `
      '<w:r><w:t>' + $closeQuote + '</w:t></w:r>' + `
      '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> or </w:t></w:r>' + `
      '<w:r><w:t>' + $openQuote + '</w:t></w:r>' + `
      '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>two raised to the power of three</w:t></w:r>' + `
      '<w:r><w:t>' + $closeQuote + '</w:t></w:r>' + `
    '</w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($xml)

# InsertXML splits the paragraph, leaving behind a duplicate (empty) copy of
# the original paragraph mark right after the freshly inserted paragraph.
# Merge it away so the document only gains the single intended paragraph.
$idx2 = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $newRange.Start) {
        $idx2 = $i
        break
    }
}
$insertedPara = $d.Paragraphs.Item($idx2)
$strayPara = $d.Paragraphs.Item($idx2 + 1)
$strayText = $strayPara.Range.Text
if ($strayText -eq "" -or $strayText -eq "`r") {
    $cleanupRange = $d.Range($insertedPara.Range.End - 1, $strayPara.Range.End)
    $cleanupRange.Delete()
}
